# Apply "All updates before first test" changes to CurrentCatheter sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Remove the "Xi (distance between pin wall to catheter wall)" column.
# Columns I1:K1 shift left by one (old J1->I1, K1->J1, L1->K1), and the
# trailing column (old L1, "Mandrel OD (mm)") is dropped.
$ws.Range("I1").Value = "Yi (Dist end of grippers to bending pin)"
$ws.Range("J1").Value = "Mandrel Material"
$ws.Range("K1").Value = "Mandrel OD (mm)"
$ws.Range("L1").Clear()

# --- Data rows ---
# Row 2 (Material A)
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "A"
$ws.Range("C2").Value = 4.25
$ws.Range("D2").Value = 1.2
$ws.Range("E2").Value = 1.98
$ws.Range("F2").Value = "Soft, black"
$ws.Range("H2").Value = 0.5
$ws.Range("I2").Value = 0.5

# Row 3 (Material B)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "B"
$ws.Range("C3").Value = 42
$ws.Range("D3").Value = 1.67
$ws.Range("E3").Value = 1.17
$ws.Range("F3").Value = "Soft, black"
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 0.5

# Row 4 (Material C)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "C"
$ws.Range("C4").Value = 100
$ws.Range("D4").Value = 1.67
$ws.Range("E4").Value = 1.17
$ws.Range("F4").Value = "Braided, purple"
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 0.5

# Row 5 (Material D) - new row
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "D"
$ws.Range("C5").Value = 50
$ws.Range("D5").Value = 1.67
$ws.Range("E5").Value = 1.17
$ws.Range("F5").Value = "Soft, black"
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 0.5

# Row 6 (Material E) - new row
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "E"
$ws.Range("C6").Value = 50
$ws.Range("D6").Value = 1.33
$ws.Range("E6").Value = 1.17
$ws.Range("F6").Value = "Soft, black"
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = 0.5

# Apply the same formatting used by A2:A4 (bold font, boxed border, centered)
# to the new A5:A6 cells by copying formats only (keeps styles.xml minimal,
# matching the style index already used by the other "A" column cells).
$ws.Range("A2").Copy()
$ws.Range("A5:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
